# Add the "grid_sellnothing_test" worksheet to the workbook.
#
# This mirrors the existing "grid_sellall_test" sheet (same component list /
# layout) but with the electricity-price bounds and expected-results values
# changed to reflect a scenario where the optimizer should sell as little
# power to the grid as possible (instead of as much as possible).

$wb = $excel.ActiveWorkbook

# --- Update the previously-active sheet's view state -----------------------
# grid_sellall_test was the selected/active tab before this edit; once the
# new sheet is added and activated it is no longer the active tab, and its
# selection is reset to the full data range.
$gridSellAll = $wb.Worksheets.Item("grid_sellall_test")
$gridSellAll.Activate() | Out-Null
$gridSellAll.Range("A1:F22").Select() | Out-Null

# --- Add the new worksheet after the last existing sheet --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "grid_sellnothing_test"
$ws = $newSheet

$ws.Range("A1").Value = "Electricity price"
$ws.Range("B1").Value = 0.1
$ws.Range("C1").Value = "`$/MWh"
$ws.Range("B2").Formula = "=B1/1000"
$ws.Range("C2").Value = "`$/kWh"
$ws.Range("A3").Value = "Naphta"
$ws.Range("B3").Value = 10000
$ws.Range("B3").NumberFormat = "0.00E+00"
$ws.Range("C3").Value = "`$/kg"
$ws.Range("A4").Value = "Jet Fuel"
$ws.Range("B4").Value = 10000
$ws.Range("B4").NumberFormat = "0.00E+00"
$ws.Range("C4").Value = "`$/kg"
$ws.Range("A5").Value = "Diesel"
$ws.Range("B5").Value = 10000
$ws.Range("B5").NumberFormat = "0.00E+00"
$ws.Range("C5").Value = "`$/kg"
$ws.Range("A7").Value = "Bounds"
$ws.Range("A8").Value = "Component"
$ws.Range("B8").Value = "Low "
$ws.Range("C8").Value = "High"
$ws.Range("D8").Value = "Unit"
$ws.Range("E8").Value = "Quantity"
$ws.Range("A9").Value = "Turbine"
$ws.Range("B9").Value = 1000
$ws.Range("B9").NumberFormat = "0.00E+00"
$ws.Range("C9").Value = 750000
$ws.Range("C9").NumberFormat = "0.00E+00"
$ws.Range("D9").Value = "kWe"
$ws.Range("E9").Value = "Elec"
$ws.Range("A10").Value = "HTSE"
$ws.Range("B10").Value = -750000
$ws.Range("B10").NumberFormat = "0.00E+00"
$ws.Range("C10").Value = -1000
$ws.Range("C10").NumberFormat = "0.00E+00"
$ws.Range("D10").Value = "kWe"
$ws.Range("E10").Value = "Elec"
$ws.Range("A11").Value = "FT"
$ws.Range("B11").Value = -100000
$ws.Range("B11").NumberFormat = "0.00E+00"
$ws.Range("C11").Value = -1000
$ws.Range("C11").NumberFormat = "0.00E+00"
$ws.Range("D11").Value = "kg/h"
$ws.Range("E11").Value = "H2"
$ws.Range("A12").Value = "CO2 source"
$ws.Range("B12").Value = 1000
$ws.Range("B12").NumberFormat = "0.00E+00"
$ws.Range("C12").Value = 500000
$ws.Range("C12").NumberFormat = "0.00E+00"
$ws.Range("D12").Value = "kg/h"
$ws.Range("E12").Value = "CO2"
$ws.Range("A13").Value = "H2 storage"
$ws.Range("B13").Value = 1000
$ws.Range("B13").NumberFormat = "0.00E+00"
$ws.Range("C13").Value = 200000
$ws.Range("C13").NumberFormat = "0.00E+00"
$ws.Range("D13").Value = "kg "
$ws.Range("E13").Value = "H2"
$ws.Range("A15").Value = "Should sell as little as possible to grid"
$ws.Range("A16").Value = "Expected results"
$ws.Range("A17").Value = "Component"
$ws.Range("B17").Value = "Optimized capacity"
$ws.Range("C17").Value = "Unit"
$ws.Range("D17").Value = "Quantity"
$ws.Range("A18").Value = "Turbine"
$ws.Range("A19").Value = "HTSE"
$ws.Range("B19").Value = -750000
$ws.Range("B19").NumberFormat = "0.00E+00"
$ws.Range("C19").Value = "kWe"
$ws.Range("D19").Value = "Elec"
$ws.Range("A20").Value = "FT"
$ws.Range("B20").Value = -100000
$ws.Range("B20").NumberFormat = "0.00E+00"
$ws.Range("C20").Value = "kg/h"
$ws.Range("D20").Value = "H2"
$ws.Range("A21").Value = "CO2 source"
$ws.Range("A22").Value = "H2 storage"
$ws.Range("B22").Value = 1000
$ws.Range("B22").NumberFormat = "0.00E+00"
$ws.Range("C22").Value = "kg "
$ws.Range("D22").Value = "H2"
$ws.Range("E22").Value = "Still smallest possible since constant prices everywhere"

# Select the last-edited cell on the new sheet, and make it the active tab
# (matching the workbook's updated activeTab index).
$newSheet.Activate() | Out-Null
$ws.Range("E22").Select() | Out-Null
